$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of old (BGR decimal) Interior.Color values to the new target colors,
# derived from the fgColor swap in xl/styles.xml <fills>.
$colorMap = @{
    6946889 = 1841892
    7209049 = 12090935
    7537002 = 4894541
    7799162 = 10702488
    7930251 = 32767
    8126877 = 3407871
    8257966 = 2643622
    8786622 = 12550647
    9380813 = 10066329
    9909469 = 10863206
    10110438 = 6458876
    10377198 = 13344909
    10578167 = 12815079
    11041528 = 5560486
    11439609 = 3135999
    11902970 = 9749733
    12168443 = 11776947
    12368123 = 7839259
    12633596 = 155609
    13291260 = 11759733
    13883389 = 9054695
    14541053 = 2008678
    15001854 = 175078
    15527934 = 1930918
    15988735 = 6710886
}

$used = $ws.UsedRange
$firstRow = $used.Row
$lastRow = $firstRow + $used.Rows.Count - 1
$firstCol = $used.Column
$lastCol = $firstCol + $used.Columns.Count - 1

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowRange = $ws.Range($ws.Cells.Item($r, $firstCol), $ws.Cells.Item($r, $lastCol))
    $current = $rowRange.Cells.Item(1, 1).Interior.Color
    if ($colorMap.ContainsKey($current)) {
        $rowRange.Interior.Color = $colorMap[$current]
    }
}
